$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("UK_cohabitation_adjustment")
$ws4 = $wb.Worksheets.Item("UK_fertility_adjustment")

$ws3.Cells.Item(1, 5).Value = "standard model"
$ws3.Cells.Item(1, 6).Value = "dynamic utility model"
$ws3.Cells.Item(2, 5).Value = 0
$ws3.Cells.Item(2, 6).Value = 0
$ws3.Cells.Item(3, 5).Value = 0
$ws3.Cells.Item(3, 6).Value = 0
$ws3.Cells.Item(4, 2).Value = -0.68042582690000331
$ws3.Cells.Item(4, 5).Value = -0.95768509819680903
$ws3.Cells.Item(4, 6).Value = -0.68042582690000331
$ws3.Cells.Item(5, 2).Value = -0.67229930191496556
$ws3.Cells.Item(5, 5).Value = -0.83295781453007001
$ws3.Cells.Item(5, 6).Value = -0.67229930191496556
$ws3.Cells.Item(6, 2).Value = -0.67897363322872639
$ws3.Cells.Item(6, 5).Value = -0.80485660959326355
$ws3.Cells.Item(6, 6).Value = -0.67897363322872639
$ws3.Cells.Item(7, 2).Value = -0.66203073162886494
$ws3.Cells.Item(7, 5).Value = -0.78981867057924415
$ws3.Cells.Item(7, 6).Value = -0.66203073162886494
$ws3.Cells.Item(8, 2).Value = -0.65127749552133463
$ws3.Cells.Item(8, 5).Value = -0.77794475503792282
$ws3.Cells.Item(8, 6).Value = -0.65127749552133463
$ws3.Cells.Item(9, 2).Value = -0.62856645997943639
$ws3.Cells.Item(9, 5).Value = -0.75171647264083308
$ws3.Cells.Item(9, 6).Value = -0.62856645997943639
$ws3.Cells.Item(10, 2).Value = -0.62262883931939683
$ws3.Cells.Item(10, 5).Value = -0.74109647439251236
$ws3.Cells.Item(10, 6).Value = -0.62262883931939683
$ws3.Cells.Item(11, 2).Value = -0.61386966551793498
$ws3.Cells.Item(11, 5).Value = -0.72776510107134573
$ws3.Cells.Item(11, 6).Value = -0.61386966551793498
$ws3.Cells.Item(12, 2).Value = -0.59884444696646777
$ws3.Cells.Item(12, 5).Value = -0.7147060419149136
$ws3.Cells.Item(12, 6).Value = -0.59884444696646777
$ws3.Cells.Item(13, 2).Value = -0.58212592838721533
$ws3.Cells.Item(13, 5).Value = -0.70240616614065376
$ws3.Cells.Item(13, 6).Value = -0.58212592838721533
$ws3.Cells.Item(14, 2).Value = -0.58476026852042595
$ws3.Cells.Item(14, 5).Value = -0.69105367590666433
$ws3.Cells.Item(14, 6).Value = -0.58476026852042595
$ws3.Cells.Item(15, 2).Value = -0.56894577426477011
$ws3.Cells.Item(15, 5).Value = -0.67898522305967135
$ws3.Cells.Item(15, 6).Value = -0.56894577426477011
$ws3.Cells.Item(16, 2).Value = -0.55774239163231498
$ws3.Cells.Item(16, 5).Value = -0.67837200227633576
$ws3.Cells.Item(16, 6).Value = -0.55774239163231498
$ws3.Cells.Item(17, 2).Value = -0.5546484315400706
$ws3.Cells.Item(17, 5).Value = -0.66865521484219381
$ws3.Cells.Item(17, 6).Value = -0.5546484315400706
$ws3.Cells.Item(18, 2).Value = -0.54633744534475526
$ws3.Cells.Item(18, 5).Value = -0.66213777125088935
$ws3.Cells.Item(18, 6).Value = -0.54633744534475526
$ws3.Cells.Item(19, 2).Value = -0.53660398326743075
$ws3.Cells.Item(19, 5).Value = -0.65624918555298262
$ws3.Cells.Item(19, 6).Value = -0.53660398326743075
$ws3.Cells.Item(20, 2).Value = -0.52636891781716488
$ws3.Cells.Item(20, 5).Value = -0.65247439548232655
$ws3.Cells.Item(20, 6).Value = -0.52636891781716488
$ws3.Cells.Item(21, 2).Value = -0.52739427357167223
$ws3.Cells.Item(21, 5).Value = -0.64776174741668557
$ws3.Cells.Item(21, 6).Value = -0.52739427357167223
$ws3.Cells.Item(22, 2).Value = -0.522335868711721
$ws3.Cells.Item(22, 5).Value = -0.64184148070049241
$ws3.Cells.Item(22, 6).Value = -0.522335868711721
$ws3.Cells.Item(23, 2).Value = -0.51932161850484515
$ws3.Cells.Item(23, 5).Value = -0.63931566276826524
$ws3.Cells.Item(23, 6).Value = -0.51932161850484515
$ws3.Cells.Item(24, 2).Value = -0.51851186880903399
$ws3.Cells.Item(24, 5).Value = -0.63638760354965007
$ws3.Cells.Item(24, 6).Value = -0.51851186880903399
$ws3.Cells.Item(25, 2).Value = -0.52248458236888451
$ws3.Cells.Item(25, 5).Value = -0.63218797948615213
$ws3.Cells.Item(25, 6).Value = -0.52248458236888451
$ws3.Cells.Item(26, 2).Value = -0.51756454907723248
$ws3.Cells.Item(26, 5).Value = -0.62979949109872213
$ws3.Cells.Item(26, 6).Value = -0.51756454907723248
$ws3.Cells.Item(27, 2).Value = -0.51707659305484055
$ws3.Cells.Item(27, 5).Value = -0.62642840586058846
$ws3.Cells.Item(27, 6).Value = -0.51707659305484055
$ws3.Cells.Item(28, 2).Value = -0.51785189045280133
$ws3.Cells.Item(28, 5).Value = -0.62703195252281929
$ws3.Cells.Item(28, 6).Value = -0.51785189045280133
$ws3.Cells.Item(29, 2).Value = -0.52133458277805855
$ws3.Cells.Item(29, 5).Value = -0.62763592693798487
$ws3.Cells.Item(29, 6).Value = -0.52133458277805855
$ws3.Cells.Item(30, 2).Value = -0.52255964231049601
$ws3.Cells.Item(30, 5).Value = -0.62901387736233438
$ws3.Cells.Item(30, 6).Value = -0.52255964231049601
$ws3.Cells.Item(31, 2).Value = -0.52461978177447943
$ws3.Cells.Item(31, 5).Value = -0.62870490298773551
$ws3.Cells.Item(31, 6).Value = -0.52461978177447943
$ws3.Cells.Item(32, 2).Value = -0.52602689591038243
$ws3.Cells.Item(32, 5).Value = -0.63132275325212173
$ws3.Cells.Item(32, 6).Value = -0.52602689591038243
$ws3.Cells.Item(33, 2).Value = -0.52818471105801446
$ws3.Cells.Item(33, 5).Value = -0.63189417814259741
$ws3.Cells.Item(33, 6).Value = -0.52818471105801446
$ws3.Cells.Item(34, 2).Value = -0.52582640158760352
$ws3.Cells.Item(34, 5).Value = -0.63431628626780423
$ws3.Cells.Item(34, 6).Value = -0.52582640158760352
$ws3.Cells.Item(35, 2).Value = -0.52884537398341225
$ws3.Cells.Item(35, 5).Value = -0.63319756552905415
$ws3.Cells.Item(35, 6).Value = -0.52884537398341225
$ws3.Cells.Item(36, 2).Value = -0.53136207561516069
$ws3.Cells.Item(36, 5).Value = -0.63471667603087323
$ws3.Cells.Item(36, 6).Value = -0.53136207561516069
$ws3.Cells.Item(37, 2).Value = -0.53392813259970762
$ws3.Cells.Item(37, 5).Value = -0.63455648067220116
$ws3.Cells.Item(37, 6).Value = -0.53392813259970762
$ws3.Cells.Item(38, 2).Value = -0.53502526357346303
$ws3.Cells.Item(38, 5).Value = -0.63527012935003857
$ws3.Cells.Item(38, 6).Value = -0.53502526357346303
$ws3.Cells.Item(39, 2).Value = -0.5361827911656164
$ws3.Cells.Item(39, 5).Value = -0.63488663198137396
$ws3.Cells.Item(39, 6).Value = -0.5361827911656164
$ws3.Cells.Item(40, 2).Value = -0.53742998980006962
$ws3.Cells.Item(40, 5).Value = -0.6327623871436846
$ws3.Cells.Item(40, 6).Value = -0.53742998980006962
$ws3.Cells.Item(41, 2).Value = -0.53265366720456808
$ws3.Cells.Item(41, 5).Value = -0.63172645653828041
$ws3.Cells.Item(41, 6).Value = -0.53265366720456808
$ws3.Cells.Item(42, 2).Value = -0.52672262067763742
$ws3.Cells.Item(42, 5).Value = -0.63131309782462763
$ws3.Cells.Item(42, 6).Value = -0.52672262067763742
$ws3.Cells.Item(43, 2).Value = -0.52246561759781496
$ws3.Cells.Item(43, 5).Value = -0.62711212192261312
$ws3.Cells.Item(43, 6).Value = -0.52246561759781496
$ws3.Cells.Item(44, 2).Value = -0.51992475724234222
$ws3.Cells.Item(44, 5).Value = -0.62483324297863929
$ws3.Cells.Item(44, 6).Value = -0.51992475724234222
$ws3.Cells.Item(45, 2).Value = -0.51600574555924328
$ws3.Cells.Item(45, 5).Value = -0.62289315139452228
$ws3.Cells.Item(45, 6).Value = -0.51600574555924328
$ws3.Cells.Item(46, 2).Value = -0.51576972424348644
$ws3.Cells.Item(46, 5).Value = -0.6208575795852076
$ws3.Cells.Item(46, 6).Value = -0.51576972424348644
$ws3.Cells.Item(47, 2).Value = -0.51800776098642443
$ws3.Cells.Item(47, 5).Value = -0.62046331508145192
$ws3.Cells.Item(47, 6).Value = -0.51800776098642443
$ws3.Cells.Item(48, 2).Value = -0.51952270196948724
$ws3.Cells.Item(48, 5).Value = -0.61860426418741976
$ws3.Cells.Item(48, 6).Value = -0.51952270196948724
$ws3.Cells.Item(49, 2).Value = -0.52097143042496585
$ws3.Cells.Item(49, 5).Value = -0.61533367788355553
$ws3.Cells.Item(49, 6).Value = -0.52097143042496585
$ws3.Cells.Item(50, 2).Value = -0.51855272309252176
$ws3.Cells.Item(50, 5).Value = -0.61489963029454253
$ws3.Cells.Item(50, 6).Value = -0.51855272309252176
$ws3.Cells.Item(51, 2).Value = -0.51809916084156005
$ws3.Cells.Item(51, 5).Value = -0.61364330813003531
$ws3.Cells.Item(51, 6).Value = -0.51809916084156005
$ws3.Cells.Item(52, 2).Value = -0.51977389607061364
$ws3.Cells.Item(52, 5).Value = -0.61242239408296961
$ws3.Cells.Item(52, 6).Value = -0.51977389607061364
$ws3.Cells.Item(53, 2).Value = -0.51599900885427696
$ws3.Cells.Item(53, 5).Value = -0.61091556383364465
$ws3.Cells.Item(53, 6).Value = -0.51599900885427696
$ws3.Cells.Item(54, 2).Value = -0.51480700328602169
$ws3.Cells.Item(54, 5).Value = -0.60932285274850906
$ws3.Cells.Item(54, 6).Value = -0.51480700328602169
$ws3.Cells.Item(55, 2).Value = -0.51568068499828201
$ws3.Cells.Item(55, 5).Value = -0.60791630071486702
$ws3.Cells.Item(55, 6).Value = -0.51568068499828201
$ws3.Cells.Item(56, 2).Value = -0.51597773825575943
$ws3.Cells.Item(56, 5).Value = -0.60902347078757146
$ws3.Cells.Item(56, 6).Value = -0.51597773825575943
$ws3.Cells.Item(57, 2).Value = -0.50918163905100955
$ws3.Cells.Item(57, 5).Value = -0.60734177454365323
$ws3.Cells.Item(57, 6).Value = -0.50918163905100955
$ws3.Cells.Item(58, 2).Value = -0.51037220329257249
$ws3.Cells.Item(58, 5).Value = -0.60672079418846891
$ws3.Cells.Item(58, 6).Value = -0.51037220329257249
$ws3.Cells.Item(59, 2).Value = -0.50578085391324401
$ws3.Cells.Item(59, 5).Value = -0.6058496118125366
$ws3.Cells.Item(59, 6).Value = -0.50578085391324401
$ws3.Cells.Item(60, 2).Value = -0.50041296740445906
$ws3.Cells.Item(60, 5).Value = -0.60652410750329455
$ws3.Cells.Item(60, 6).Value = -0.50041296740445906
$ws3.Cells.Item(61, 2).Value = -0.50041296740445906
$ws3.Cells.Item(61, 5).Value = -0.60652410750329455
$ws3.Cells.Item(61, 6).Value = -0.50041296740445906
$ws3.Cells.Item(62, 2).Value = -0.50041296740445906
$ws3.Cells.Item(62, 5).Value = -0.60652410750329455
$ws3.Cells.Item(62, 6).Value = -0.50041296740445906

$ws4.Cells.Item(1, 5).Value = "standard model"
$ws4.Cells.Item(1, 6).Value = "dynamic utility model"
$ws4.Cells.Item(2, 5).Value = 0
$ws4.Cells.Item(2, 6).Value = 0
$ws4.Cells.Item(3, 5).Value = 0
$ws4.Cells.Item(3, 6).Value = 0
$ws4.Cells.Item(4, 2).Value = -0.44081411312763502
$ws4.Cells.Item(4, 5).Value = -0.48463199296432402
$ws4.Cells.Item(4, 6).Value = -0.44081411312763502
$ws4.Cells.Item(5, 2).Value = -0.41178003442893141
$ws4.Cells.Item(5, 5).Value = -0.44279476957698699
$ws4.Cells.Item(5, 6).Value = -0.41178003442893141
$ws4.Cells.Item(6, 2).Value = -0.38195691653651342
$ws4.Cells.Item(6, 5).Value = -0.39501557684159466
$ws4.Cells.Item(6, 6).Value = -0.38195691653651342
$ws4.Cells.Item(7, 2).Value = -0.35662154819469444
$ws4.Cells.Item(7, 5).Value = -0.37898938142570604
$ws4.Cells.Item(7, 6).Value = -0.35662154819469444
$ws4.Cells.Item(8, 2).Value = -0.33645665993640061
$ws4.Cells.Item(8, 5).Value = -0.35746845133604671
$ws4.Cells.Item(8, 6).Value = -0.33645665993640061
$ws4.Cells.Item(9, 2).Value = -0.31569377810596422
$ws4.Cells.Item(9, 5).Value = -0.33551640350223444
$ws4.Cells.Item(9, 6).Value = -0.31569377810596422
$ws4.Cells.Item(10, 2).Value = -0.29611480221652442
$ws4.Cells.Item(10, 5).Value = -0.31611792786789245
$ws4.Cells.Item(10, 6).Value = -0.29611480221652442
$ws4.Cells.Item(11, 2).Value = -0.279462622342386
$ws4.Cells.Item(11, 5).Value = -0.29419113555012588
$ws4.Cells.Item(11, 6).Value = -0.279462622342386
$ws4.Cells.Item(12, 2).Value = -0.26591979180000819
$ws4.Cells.Item(12, 5).Value = -0.275174972273524
$ws4.Cells.Item(12, 6).Value = -0.26591979180000819
$ws4.Cells.Item(13, 2).Value = -0.25543875103684177
$ws4.Cells.Item(13, 5).Value = -0.25957357188884783
$ws4.Cells.Item(13, 6).Value = -0.25543875103684177
$ws4.Cells.Item(14, 2).Value = -0.24930347300809638
$ws4.Cells.Item(14, 5).Value = -0.2522555041315423
$ws4.Cells.Item(14, 6).Value = -0.24930347300809638
$ws4.Cells.Item(15, 2).Value = -0.2471951557625644
$ws4.Cells.Item(15, 5).Value = -0.2484037751560034
$ws4.Cells.Item(15, 6).Value = -0.2471951557625644
$ws4.Cells.Item(16, 2).Value = -0.24727245707738721
$ws4.Cells.Item(16, 5).Value = -0.24178593173429327
$ws4.Cells.Item(16, 6).Value = -0.24727245707738721
$ws4.Cells.Item(17, 2).Value = -0.23666450179745918
$ws4.Cells.Item(17, 5).Value = -0.24449858212979714
$ws4.Cells.Item(17, 6).Value = -0.23666450179745918
$ws4.Cells.Item(18, 2).Value = -0.23784361977757321
$ws4.Cells.Item(18, 5).Value = -0.24903198566118756
$ws4.Cells.Item(18, 6).Value = -0.23784361977757321
$ws4.Cells.Item(19, 2).Value = -0.23889278870939884
$ws4.Cells.Item(19, 5).Value = -0.25264201070565656
$ws4.Cells.Item(19, 6).Value = -0.23889278870939884
$ws4.Cells.Item(20, 2).Value = -0.23893998110579662
$ws4.Cells.Item(20, 5).Value = -0.25364233752894055
$ws4.Cells.Item(20, 6).Value = -0.23893998110579662
$ws4.Cells.Item(21, 2).Value = -0.23843847695520243
$ws4.Cells.Item(21, 5).Value = -0.24768152163486329
$ws4.Cells.Item(21, 6).Value = -0.23843847695520243
$ws4.Cells.Item(22, 2).Value = -0.24968508459850719
$ws4.Cells.Item(22, 5).Value = -0.24974442678892514
$ws4.Cells.Item(22, 6).Value = -0.24968508459850719
$ws4.Cells.Item(23, 2).Value = -0.25091556848125485
$ws4.Cells.Item(23, 5).Value = -0.24926794764317603
$ws4.Cells.Item(23, 6).Value = -0.25091556848125485
$ws4.Cells.Item(24, 2).Value = -0.25300049353502324
$ws4.Cells.Item(24, 5).Value = -0.250134489952685
$ws4.Cells.Item(24, 6).Value = -0.25300049353502324
$ws4.Cells.Item(25, 2).Value = -0.25769551121458001
$ws4.Cells.Item(25, 5).Value = -0.25159703697288732
$ws4.Cells.Item(25, 6).Value = -0.25769551121458001
$ws4.Cells.Item(26, 2).Value = -0.26297059722398985
$ws4.Cells.Item(26, 5).Value = -0.25383394104860174
$ws4.Cells.Item(26, 6).Value = -0.26297059722398985
$ws4.Cells.Item(27, 2).Value = -0.26957256055242024
$ws4.Cells.Item(27, 5).Value = -0.26016905205776741
$ws4.Cells.Item(27, 6).Value = -0.26957256055242024
$ws4.Cells.Item(28, 2).Value = -0.27813755313923438
$ws4.Cells.Item(28, 5).Value = -0.27274303604096028
$ws4.Cells.Item(28, 6).Value = -0.27813755313923438
$ws4.Cells.Item(29, 2).Value = -0.27474637053283057
$ws4.Cells.Item(29, 5).Value = -0.27650991999973468
$ws4.Cells.Item(29, 6).Value = -0.27474637053283057
$ws4.Cells.Item(30, 2).Value = -0.28518071039366039
$ws4.Cells.Item(30, 5).Value = -0.28460685361503957
$ws4.Cells.Item(30, 6).Value = -0.28518071039366039
$ws4.Cells.Item(31, 2).Value = -0.2959050156529302
$ws4.Cells.Item(31, 5).Value = -0.29665890398267214
$ws4.Cells.Item(31, 6).Value = -0.2959050156529302
$ws4.Cells.Item(32, 2).Value = -0.2928526070865794
$ws4.Cells.Item(32, 5).Value = -0.30980188741828268
$ws4.Cells.Item(32, 6).Value = -0.2928526070865794
$ws4.Cells.Item(33, 2).Value = -0.30204917920681157
$ws4.Cells.Item(33, 5).Value = -0.3249147353451094
$ws4.Cells.Item(33, 6).Value = -0.30204917920681157
$ws4.Cells.Item(34, 2).Value = -0.32304908725990761
$ws4.Cells.Item(34, 5).Value = -0.33013144629608343
$ws4.Cells.Item(34, 6).Value = -0.32304908725990761
$ws4.Cells.Item(35, 2).Value = -0.32697846249115198
$ws4.Cells.Item(35, 5).Value = -0.33400456174430643
$ws4.Cells.Item(35, 6).Value = -0.32697846249115198
$ws4.Cells.Item(36, 2).Value = -0.33007523822817458
$ws4.Cells.Item(36, 5).Value = -0.34264176458383305
$ws4.Cells.Item(36, 6).Value = -0.33007523822817458
$ws4.Cells.Item(37, 2).Value = -0.34395374911145921
$ws4.Cells.Item(37, 5).Value = -0.34821267929122846
$ws4.Cells.Item(37, 6).Value = -0.34395374911145921
$ws4.Cells.Item(38, 2).Value = -0.32911187304579681
$ws4.Cells.Item(38, 5).Value = -0.34456161150133158
$ws4.Cells.Item(38, 6).Value = -0.32911187304579681
$ws4.Cells.Item(39, 2).Value = -0.32707905887567962
$ws4.Cells.Item(39, 5).Value = -0.34561532498119613
$ws4.Cells.Item(39, 6).Value = -0.32707905887567962
$ws4.Cells.Item(40, 2).Value = -0.32548082652979621
$ws4.Cells.Item(40, 5).Value = -0.33873736519970848
$ws4.Cells.Item(40, 6).Value = -0.32548082652979621
$ws4.Cells.Item(41, 2).Value = -0.32367111158351802
$ws4.Cells.Item(41, 5).Value = -0.33739651765499096
$ws4.Cells.Item(41, 6).Value = -0.32367111158351802
$ws4.Cells.Item(42, 2).Value = -0.31183478961522498
$ws4.Cells.Item(42, 5).Value = -0.33804860844064244
$ws4.Cells.Item(42, 6).Value = -0.31183478961522498
$ws4.Cells.Item(43, 2).Value = -0.30933783525634517
$ws4.Cells.Item(43, 5).Value = -0.33657044977191414
$ws4.Cells.Item(43, 6).Value = -0.30933783525634517
$ws4.Cells.Item(44, 2).Value = -0.295502572399851
$ws4.Cells.Item(44, 5).Value = -0.33656675024980498
$ws4.Cells.Item(44, 6).Value = -0.295502572399851
$ws4.Cells.Item(45, 2).Value = -0.28426566851705959
$ws4.Cells.Item(45, 5).Value = -0.33985483078933226
$ws4.Cells.Item(45, 6).Value = -0.28426566851705959
$ws4.Cells.Item(46, 2).Value = -0.2846970607117626
$ws4.Cells.Item(46, 5).Value = -0.33668670713123178
$ws4.Cells.Item(46, 6).Value = -0.2846970607117626
$ws4.Cells.Item(47, 2).Value = -0.29652418423720162
$ws4.Cells.Item(47, 5).Value = -0.33750504731956005
$ws4.Cells.Item(47, 6).Value = -0.29652418423720162
$ws4.Cells.Item(48, 2).Value = -0.31354110133198604
$ws4.Cells.Item(48, 5).Value = -0.34079779337557842
$ws4.Cells.Item(48, 6).Value = -0.31354110133198604
$ws4.Cells.Item(49, 2).Value = -0.3278678462801522
$ws4.Cells.Item(49, 5).Value = -0.33872371122332023
$ws4.Cells.Item(49, 6).Value = -0.3278678462801522
$ws4.Cells.Item(50, 2).Value = -0.33932540739172184
$ws4.Cells.Item(50, 5).Value = -0.33869137298543883
$ws4.Cells.Item(50, 6).Value = -0.33932540739172184
$ws4.Cells.Item(51, 2).Value = -0.33891825783445723
$ws4.Cells.Item(51, 5).Value = -0.34103219793504497
$ws4.Cells.Item(51, 6).Value = -0.33891825783445723
$ws4.Cells.Item(52, 2).Value = -0.33791821188593724
$ws4.Cells.Item(52, 5).Value = -0.34050601889900228
$ws4.Cells.Item(52, 6).Value = -0.33791821188593724
$ws4.Cells.Item(53, 2).Value = -0.33627771003385082
$ws4.Cells.Item(53, 5).Value = -0.33767373591580502
$ws4.Cells.Item(53, 6).Value = -0.33627771003385082
$ws4.Cells.Item(54, 2).Value = -0.33352763519739703
$ws4.Cells.Item(54, 5).Value = -0.33668796343699775
$ws4.Cells.Item(54, 6).Value = -0.33352763519739703
$ws4.Cells.Item(55, 2).Value = -0.33028054032751103
$ws4.Cells.Item(55, 5).Value = -0.33369114324283816
$ws4.Cells.Item(55, 6).Value = -0.33028054032751103
$ws4.Cells.Item(56, 2).Value = -0.32636589095194396
$ws4.Cells.Item(56, 5).Value = -0.33282968867461227
$ws4.Cells.Item(56, 6).Value = -0.32636589095194396
$ws4.Cells.Item(57, 2).Value = -0.32296966045752384
$ws4.Cells.Item(57, 5).Value = -0.32675564471673357
$ws4.Cells.Item(57, 6).Value = -0.32296966045752384
$ws4.Cells.Item(58, 2).Value = -0.3197457208767488
$ws4.Cells.Item(58, 5).Value = -0.32143826156637562
$ws4.Cells.Item(58, 6).Value = -0.3197457208767488
$ws4.Cells.Item(59, 2).Value = -0.31773659008754962
$ws4.Cells.Item(59, 5).Value = -0.31711471655916201
$ws4.Cells.Item(59, 6).Value = -0.31773659008754962
$ws4.Cells.Item(60, 2).Value = -0.31670999077073703
$ws4.Cells.Item(60, 5).Value = -0.31670999077073719
$ws4.Cells.Item(60, 6).Value = -0.31670999077073703
$ws4.Cells.Item(61, 2).Value = -0.31670999077073703
$ws4.Cells.Item(61, 5).Value = -0.31670999077073719
$ws4.Cells.Item(61, 6).Value = -0.31670999077073703
$ws4.Cells.Item(62, 2).Value = -0.31670999077073703
$ws4.Cells.Item(62, 5).Value = -0.31670999077073719
$ws4.Cells.Item(62, 6).Value = -0.31670999077073703

# Re-create the per-sheet selection state (B2:B62) while preserving which
# sheet is the active/selected tab (UK_fertility_adjustment, same as before).
[void]$ws3.Select()
[void]$ws3.Range("B2:B62").Select()
[void]$ws4.Select()
[void]$ws4.Range("B2:B62").Select()
